# Apply the "changed general parameter to dictionary" edit.
$wb = $excel.ActiveWorkbook

# 1. Remove the unused helper sheets (Tabelle1/2/3).
$wb.Worksheets("Tabelle2").Delete() | Out-Null
$wb.Worksheets("Tabelle1").Delete() | Out-Null
$wb.Worksheets("Tabelle3").Delete() | Out-Null

# 2. Rename "Costs of default system" -> "Costs default system".
$wb.Worksheets("Costs of default system").Name = "Costs default system"

# 3. "Sets" sheet: add a new "Cost type default" column (F), mirroring the
#    "Cost type" column (E) but shifted up by one row.
$wsSets = $wb.Worksheets("Sets")
$wsSets.Range("F1").Value = "Cost type default"
$wsSets.Range("F2").Value = "Service Cost"
$wsSets.Range("F3").Value = "Connection Price"
$wsSets.Range("F4").Value = "Fuel Price"
$wsSets.Range("F5").Value = "Feedin Price"

# 4. "General Data" sheet: fix a label typo and repurpose the charging
#    station count row into a car count row.
$wsGeneral = $wb.Worksheets("General Data")
$wsGeneral.Range("A6").Value = "DHW p.P."
$wsGeneral.Range("A15").Value = "Number of cars"
$wsGeneral.Range("B15").Value = "n_cars_total"

# 5. "Costs default system" sheet: drop the Investment Price column (B) -
#    a default system is never "invested in", so the column was always 0.
$wsCostsDefault = $wb.Worksheets("Costs default system")
$wsCostsDefault.Columns("B").Delete() | Out-Null

# 6. Update selections to match the saved state of each sheet.
$wsSets.Range("D2:D4").Select()
$wsGeneral.Range("C18").Select()
$wsCostsDefault.Range("E13").Select()

# 7. Make "Demand" the active sheet/tab (was "Sets").
$wb.Worksheets("Demand").Activate()
